{"js": "// \"psychological rather than psychology\"\n//\n// 1) \"Department of Psychology Sciences\" -> \"Department of Psychological Sciences\"\n// 2) Cosmetic run-merge in the \"Publications\" hyperlink: the three runs\n//    \" google scholar\" + \" \" + \"publication list.\" collapse into a single run\n//    \" google scholar publication list.\" (visible text is unchanged).\n\nconst body = context.document.body;\n\n// --- Change 1 -------------------------------------------------------------\nconst heading = body.search(\"Department of Psychology Sciences\", { matchCase: true }).getFirst();\nheading.load(\"text\");\nawait context.sync();\n\nheading.insertText(\"Department of Psychological Sciences\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Change 2 -------------------------------------------------------------\n// The hyperlink text \" google scholar\" / \" \" / \"publication list.\" is split across three\n// runs even though the rendered text is unchanged. Assigning identical text is a no-op, so\n// round-trip the last run through a disjoint placeholder to force a real edit; the engine\n// then collapses the adjacent, identically-formatted runs into one (the preceding \"my\" run\n// stays separate because a <w:proofErr/> sits between it and the hyperlink text runs).\nconst pubRun = body.search(\"publication list.\", { matchCase: true }).getFirst();\npubRun.load(\"text\");\nawait context.sync();\n\npubRun.insertText(\"ZZZMARKERZZZ\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst marker = body.search(\"ZZZMARKERZZZ\", { matchCase: true }).getFirst();\nmarker.load(\"text\");\nawait context.sync();\n\nmarker.insertText(\"publication list.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: \"Department of Psychology Sciences\" -> \"Department of Psychological Sciences\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.MatchWildcards = $false\n$rng1.Find.Text = \"Department of Psychology Sciences\"\n$rng1.Find.Execute() | Out-Null\n$rng1.Text = \"Department of Psychological Sciences\"\n\n# --- Change 2: merge the \"publication list.\" run with its preceding sibling runs inside the\n# hyperlink (\" google scholar\" + \" \" + \"publication list.\" -> \" google scholar publication list.\")\n# The visible text does not change, so assigning the identical text is a no-op; force a real\n# edit by round-tripping through a temporary marker, which collapses the adjacent\n# identically-formatted runs into one.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.MatchWildcards = $false\n$rng2.Find.Text = \"publication list.\"\n$rng2.Find.Execute() | Out-Null\n$rng2.Text = \"publication list.ZZZMARKERZZZ\"\n\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.MatchWildcards = $false\n$rng3.Find.Text = \"publication list.ZZZMARKERZZZ\"\n$rng3.Find.Execute() | Out-Null\n$rng3.Text = \"publication list.\"\n"}
